$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@(2565, 10, "10x10", 19, 2, 9, 7, 1, 10.52631578947368, 47.36842105263158, 36.84210526315789, 5.263157894736842, "Win", "Cluster Inference")
    ,@(2566, 10, "10x10", 17, 0, 8, 8, 1, 0, 47.05882352941176, 47.05882352941176, 5.88235294117647, "Win", "Cluster Inference")
    ,@(2567, 10, "10x10", 35, 0, 28, 6, 1, 0, 80, 17.14285714285714, 2.857142857142857, "Win", "Cluster Inference")
    ,@(2568, 10, "10x10", 21, 7, 7, 5, 2, 33.33333333333333, 33.33333333333333, 23.80952380952381, 9.523809523809524, "Loss", "Monte Carlo")
    ,@(2569, 10, "10x10", 36, 0, 35, 0, 1, 0, 97.22222222222221, 0, 2.777777777777778, "Win", "Cluster Inference")
    ,@(2570, 10, "10x10", 42, 6, 26, 9, 1, 14.28571428571428, 61.90476190476191, 21.42857142857143, 2.380952380952381, "Win", "Cluster Inference")
    ,@(2571, 10, "10x10", 28, 4, 14, 9, 1, 14.28571428571428, 50, 32.14285714285715, 3.571428571428571, "Win", "Cluster Inference")
    ,@(2572, 10, "10x10", 41, 2, 32, 6, 1, 4.878048780487805, 78.04878048780488, 14.63414634146341, 2.439024390243902, "Win", "Neighbour Deduction")
    ,@(2573, 10, "10x10", 1, 0, 0, 0, 1, 0, 0, 0, 100, "Loss", "Monte Carlo")
    ,@(2574, 10, "10x10", 23, 0, 22, 0, 1, 0, 95.65217391304348, 0, 4.347826086956522, "Win", "Cluster Inference")
    ,@(2575, 10, "10x10", 33, 2, 20, 9, 2, 6.060606060606061, 60.60606060606061, 27.27272727272727, 6.060606060606061, "Win", "Cluster Inference")
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $rowNum = 119 + $i
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $rowVals.Count; $j++) {
        $colNum = $j + 1
        $ws.Cells.Item($rowNum, $colNum).Value = $rowVals[$j]
    }
}
